# Scheduled-runner market data refresh for Kujata_Profits workbook.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) per sheet/row,
# clearing cells that no longer carry a computed value (ClearContents).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3875
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 3875
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H107").Value = 2586.6
$ws.Range("I107").Value = 2258.25
$ws.Range("J107").Value = 3900
$ws.Range("K107").Value = 2258.25
$ws.Range("L107").Value = 3900
$ws.Range("M107").Value = -338.25
$ws.Range("N107").Value = -7740

$ws.Range("H116").Value = 3714
$ws.Range("I116").Value = 4033.3333
$ws.Range("J116").Value = 3235
$ws.Range("K116").Value = 4033.3333
$ws.Range("L116").Value = 3235
$ws.Range("M116").Value = -591.3332999999998
$ws.Range("N116").Value = -10119

$ws.Range("H137").Value = 1236.081
$ws.Range("I137").Value = 885.7
$ws.Range("J137").Value = 1648.2941
$ws.Range("K137").Value = 2657.1
$ws.Range("L137").Value = 4944.8823
$ws.Range("M137").Value = -107.1000000000004
$ws.Range("N137").Value = -10044.8823

$ws.Range("H138").Value = 507276.88
$ws.Range("I138").Value = 1417.4286
$ws.Range("J138").Value = 751484.9
$ws.Range("K138").Value = 4252.2858
$ws.Range("L138").Value = 2254454.7
$ws.Range("M138").Value = 887.7142000000003
$ws.Range("N138").Value = -2264734.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 55004.332
$ws.Range("I23").Value = 42503
$ws.Range("K23").Value = 42503
$ws.Range("M23").Value = -42244

$ws.Range("H32").Value = 3312.1924
$ws.Range("I32").Value = 2975.5217
$ws.Range("K32").Value = 2975.5217
$ws.Range("M32").Value = -2688.5217

$ws.Range("H61").Value = 1371.3939
$ws.Range("I61").Value = 1188.4286
$ws.Range("K61").Value = 1188.4286
$ws.Range("M61").Value = -976.4286

$ws.Range("H74").Value = 1408.3182
$ws.Range("I74").Value = 865.0625
$ws.Range("K74").Value = 865.0625
$ws.Range("M74").Value = 8.9375

$ws.Range("H77").Value = 1408.3182
$ws.Range("I77").Value = 865.0625
$ws.Range("K77").Value = 4325.3125
$ws.Range("M77").Value = 42.6875

$ws.Range("H110").Value = 1375.1482
$ws.Range("I110").Value = 972.17645
$ws.Range("K110").Value = 972.17645
$ws.Range("M110").Value = 1072.82355

$ws.Range("H136").Value = 1371.3939
$ws.Range("I136").Value = 1188.4286
$ws.Range("K136").Value = 3565.2858
$ws.Range("M136").Value = -1015.2858

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2950.3333
$ws.Range("I86").Value = 3138.96
$ws.Range("J86").Value = 2360.875
$ws.Range("K86").Value = 3138.96
$ws.Range("L86").Value = 2360.875
$ws.Range("M86").Value = -2015.96
$ws.Range("N86").Value = -4606.875

$ws.Range("H89").Value = 2950.3333
$ws.Range("I89").Value = 3138.96
$ws.Range("J89").Value = 2360.875
$ws.Range("K89").Value = 15694.8
$ws.Range("L89").Value = 11804.375
$ws.Range("M89").Value = -10078.8
$ws.Range("N89").Value = -23036.375

$ws.Range("H105").Value = 333336400
$ws.Range("I105").Value = 333336400
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 333336400
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -333334653
$ws.Range("N105").ClearContents()

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1262.1364
$ws.Range("I31").Value = 1196.1154
$ws.Range("J31").Value = 1507.3572
$ws.Range("K31").Value = 1196.1154
$ws.Range("L31").Value = 1507.3572
$ws.Range("M31").Value = -901.1153999999999
$ws.Range("N31").Value = -2097.3572

$ws.Range("H34").Value = 1262.1364
$ws.Range("I34").Value = 1196.1154
$ws.Range("J34").Value = 1507.3572
$ws.Range("K34").Value = 1196.1154
$ws.Range("L34").Value = 1507.3572
$ws.Range("M34").Value = -994.1153999999999
$ws.Range("N34").Value = -1911.3572

$ws.Range("H58").Value = 809.69446
$ws.Range("I58").Value = 730.6129
$ws.Range("J58").Value = 1300
$ws.Range("K58").Value = 730.6129
$ws.Range("L58").Value = 1300
$ws.Range("M58").Value = -527.6129
$ws.Range("N58").Value = -1706

$ws.Range("H136").Value = 809.69446
$ws.Range("I136").Value = 730.6129
$ws.Range("J136").Value = 1300
$ws.Range("K136").Value = 2191.8387
$ws.Range("L136").Value = 3900
$ws.Range("M136").Value = 358.1613000000002
$ws.Range("N136").Value = -9000

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 2183.2727
$ws.Range("I60").Value = 670
$ws.Range("K60").Value = 2010
$ws.Range("M60").Value = -1759

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

$ws.Range("H103").Value = 33499.668
$ws.Range("J103").Value = 33499.668
$ws.Range("L103").Value = 33499.668
$ws.Range("N103").Value = -35843.668

$ws.Range("H113").Value = 1406.5
$ws.Range("I113").Value = 1302.9
$ws.Range("J113").Value = 1579.1666
$ws.Range("K113").Value = 1302.9
$ws.Range("L113").Value = 1579.1666
$ws.Range("M113").Value = 867.0999999999999
$ws.Range("N113").Value = -5919.1666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1777.1428
$ws.Range("I68").Value = 1492
$ws.Range("J68").Value = 2822.6667
$ws.Range("K68").Value = 1492
$ws.Range("L68").Value = 2822.6667
$ws.Range("M68").Value = -743
$ws.Range("N68").Value = -4320.6667

$ws.Range("H71").Value = 1777.1428
$ws.Range("I71").Value = 1492
$ws.Range("J71").Value = 2822.6667
$ws.Range("K71").Value = 7460
$ws.Range("L71").Value = 14113.3335
$ws.Range("M71").Value = -3716
$ws.Range("N71").Value = -21601.3335

$ws.Range("H82").Value = 2285
$ws.Range("I82").Value = 2414
$ws.Range("J82").Value = 2156
$ws.Range("K82").Value = 2414
$ws.Range("L82").Value = 2156
$ws.Range("M82").Value = -2053
$ws.Range("N82").Value = -2878

$ws.Range("H85").Value = 2285
$ws.Range("I85").Value = 2414
$ws.Range("J85").Value = 2156
$ws.Range("K85").Value = 2414
$ws.Range("L85").Value = 2156
$ws.Range("M85").Value = -1166
$ws.Range("N85").Value = -4652

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 38466370
$ws.Range("I62").Value = 100006150
$ws.Range("K62").Value = 100006150
$ws.Range("M62").Value = -100005526

$ws.Range("H65").Value = 38466370
$ws.Range("I65").Value = 100006150
$ws.Range("K65").Value = 500030750
$ws.Range("M65").Value = -500027630

$ws.Range("H132").Value = 3254.9143
$ws.Range("I132").Value = 3622.0435
$ws.Range("K132").Value = 10866.1305
$ws.Range("M132").Value = -8336.130500000001
